$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.552.80"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "2.412.69"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "564.23"
$ws.Range("E5").Value = "  -3.26%  "
$ws.Range("D6").Value = "137.41"
$ws.Range("E6").Value = "  -3.66%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "2.397.18"
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("E10").Value = "  -5.51%  "
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("D14").Value = "25.65"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D16").Value = "0.0000166"
$ws.Range("E16").Value = "  -4.13%  "
$ws.Range("D17").Value = "60.611.41"
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("D18").Value = "2.380.27"
$ws.Range("E18").Value = "  -3.15%  "
$ws.Range("D19").Value = "8.16"
$ws.Range("E19").Value = "  +11.26%  "
$ws.Range("D20").Value = "10.53"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").Value = "322.12"
$ws.Range("E21").Value = "  -1.47%  "
$ws.Range("D22").Value = "4.03"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("E23").Value = "  -1.70%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  -8.31%  "
$ws.Range("D26").Value = "64.14"
$ws.Range("E26").Value = "  -1.92%  "
$ws.Range("D27").Value = "551.17"
$ws.Range("E27").Value = "  -5.69%  "
$ws.Range("D28").Value = "8.04"
$ws.Range("E28").Value = "  -12.52%  "
$ws.Range("D29").Value = "2.526.64"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").Value = "0.0₃0906"
$ws.Range("E30").Value = "  -4.19%  "
$ws.Range("D31").Value = "7.86"
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("E32").Value = "  -7.05%  "
$ws.Range("E33").Value = "  -5.18%  "
$ws.Range("E34").Value = "  -2.93%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D37").Value = "152.54"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").Value = "0.367"
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("D39").Value = "4.49"
$ws.Range("E39").Value = "  -6.24%  "
$ws.Range("D40").Value = "18.18"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("E41").Value = "  -3.03%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -4.63%  "
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "0.0₆0288"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.27"
$ws.Range("E45").Value = "  -5.61%  "
$ws.Range("D46").Value = "141.86"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").Value = "3.48"
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").Value = "0.582"
$ws.Range("E48").Value = "  -3.56%  "
$ws.Range("D50").Value = "18.98"
$ws.Range("E50").Value = "  -4.87%  "
$ws.Range("E51").Value = "  -0.84%  "
